$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "   titleS  "
$ws.Range("B5").Select()
